$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches original file, which stores all data values as shared strings)
$ws.Range("A2:H3").NumberFormat = "@"

# Update existing row 2 with the new prediction values
$ws.Range("A2").Value = "39.1"
$ws.Range("B2").Value = "181.0"
$ws.Range("C2").Value = "Adelie"
$ws.Range("D2").Value = "1.0"
$ws.Range("E2").Value = "0.0"
$ws.Range("F2").Value = "0.0"
$ws.Range("G2").Value = "v1.0"
$ws.Range("H2").Value = "2025-05-04 20:32:52"

# Add new row 3 with a new prediction
$ws.Range("A3").Value = "46.5"
$ws.Range("B3").Value = "192.0"
$ws.Range("C3").Value = "Chinstrap"
$ws.Range("D3").Value = "0.09"
$ws.Range("E3").Value = "0.91"
$ws.Range("F3").Value = "0.0"
$ws.Range("G3").Value = "v1.0"
$ws.Range("H3").Value = "2025-05-04 20:32:52"
